$d = $word.ActiveDocument

# On the title page, the subtitle paragraph reads:
#   "Syntax of SASL Symbol  Name Mangling V1.1"
# where "V1.1" is held in two separate trailing runs ("V1." and "1")
# appended after a run holding just a trailing space. The commit removes
# those two "version tag" runs entirely, leaving the paragraph ending
# right after that space run. "V1.1" occurs exactly once in the whole
# document, so a simple whole-document Find locates it unambiguously;
# deleting the matched range removes the runs that produced it without
# touching the preceding space run.
$rng = $d.Content
$found = $rng.Find.Execute("V1.1", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if ($found) {
    $rng.Delete()
}
